$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Hourly crypto-price/volume refresh (GitHub Actions bot).
# D/E columns are plain text cells (price strings use "." as a thousands
# separator, e.g. "67.880.85", and volume-change strings keep padding
# spaces, e.g. "  -0.68%  "), so every write below is a literal string.

$ws.Range("D2").Value = "67.880.85"
$ws.Range("E2").Value = "  -0.68%  "
$ws.Range("D3").Value = "3.741.56"
$ws.Range("E3").Value = "  -2.75%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "594.40"
$ws.Range("E5").Value = "  -1.22%  "
$ws.Range("D6").Value = "165.99"
$ws.Range("E6").Value = "  -3.61%  "
$ws.Range("D7").Value = "3.744.16"
$ws.Range("E7").Value = "  -2.70%  "
$ws.Range("E8").Value = "  -0.11%  "
$ws.Range("E9").Value = "  -1.91%  "
$ws.Range("E10").Value = "  -4.41%  "
$ws.Range("D11").Value = "6.44"
$ws.Range("E11").Value = "  -1.55%  "
$ws.Range("E12").Value = "  -3.38%  "
$ws.Range("E13").Value = "  -6.44%  "
$ws.Range("D14").Value = "35.95"
$ws.Range("E14").Value = "  -3.06%  "
$ws.Range("D15").Value = "4.370.40"
$ws.Range("E15").Value = "  -2.72%  "
$ws.Range("D16").Value = "3.742.41"
$ws.Range("E16").Value = "  -2.30%  "
$ws.Range("D17").Value = "67.834.59"
$ws.Range("E17").Value = "  -0.82%  "
$ws.Range("D18").Value = "18.28"
$ws.Range("E18").Value = "  -0.48%  "
$ws.Range("D19").Value = "7.05"
$ws.Range("E19").Value = "  -5.53%  "
$ws.Range("E20").Value = "  -0.15%  "
$ws.Range("D21").Value = "10.52"
$ws.Range("E21").Value = "  -3.66%  "
$ws.Range("D22").Value = "465.02"
$ws.Range("E22").Value = "  -1.32%  "
$ws.Range("E23").Value = "  -4.25%  "
$ws.Range("D24").Value = "82.92"
$ws.Range("E24").Value = "  -0.97%  "
$ws.Range("D25").Value = "0.0000138"
$ws.Range("E25").Value = "  -13.32%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.20"
$ws.Range("E26").Value = "  -3.77%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "12.00"
$ws.Range("E27").Value = "  -1.59%  "
$ws.Range("D28").Value = "10.27"
$ws.Range("E28").Value = "  -2.70%  "
$ws.Range("E29").Value = "  +0.03%  "
$ws.Range("E30").Value = "  -1.87%  "
$ws.Range("D31").Value = "3.888.26"
$ws.Range("E31").Value = "  -2.79%  "
$ws.Range("E32").Value = "  -4.93%  "
$ws.Range("D33").Value = "29.96"
$ws.Range("E33").Value = "  -3.93%  "
$ws.Range("D34").Value = "2.19"
$ws.Range("E34").Value = "  -5.14%  "
$ws.Range("D35").Value = "9.08"
$ws.Range("E35").Value = "  -3.48%  "
$ws.Range("D36").Value = "3.691.32"
$ws.Range("E36").Value = "  -3.14%  "
$ws.Range("E37").Value = "  -3.19%  "
$ws.Range("D38").Value = "3.52"
$ws.Range("E38").Value = "  -10.46%  "
$ws.Range("E39").Value = "  -1.69%  "
$ws.Range("D40").Value = "0.995"
$ws.Range("E40").Value = "  -2.37%  "
$ws.Range("D41").Value = "5.75"
$ws.Range("E41").Value = "  -3.80%  "
$ws.Range("E43").Value = "  +0.03%  "
$ws.Range("D44").Value = "0.308"
$ws.Range("E44").Value = "  -3.76%  "
$ws.Range("D45").Value = "8.54"
$ws.Range("E45").Value = "  -2.54%  "
$ws.Range("E46").Value = "  -3.79%  "
$ws.Range("D47").Value = "397.09"
$ws.Range("E47").Value = "  -5.19%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "45.00"
$ws.Range("E48").Value = "  -3.77%  "
$ws.Range("D49").Value = "143.88"
$ws.Range("E49").Value = "  +1.28%  "
$ws.Range("D50").Value = "38.87"
$ws.Range("E50").Value = "  +1.33%  "
$ws.Range("E51").Value = "  -3.64%  "
